$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3 to reference the new AMI entity's PK
$ws.Range("D3").Value = "AMI"

# Add the new row representing the AMI entity
$ws.Range("A5").Value = "Amazon Machine Image"
$ws.Range("B5").Value = "AMI"
$ws.Range("C5").Value = "AMI ID"

# Update E3 to reference the new AMI entity's SK
$ws.Range("E3").Value = "AMI ID"

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:G5"))

# Update the active selection (cosmetic change matching the diff)
$ws.Range("D8").Select()
